$d = $word.ActiveDocument

# --- Title date line -------------------------------------------------
$d.Content.Find.Execute("2025-11-25 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2025-11-26 Wednesday", 2) | Out-Null

# --- Division problems table ------------------------------------------
# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 contain
# problems (the other rows are blank spacer rows). Addressing cells
# directly (instead of a global find/replace) avoids any ambiguity from
# the repeated "685÷9=" text that must become two different values.
$t = $d.Tables.Item(1)

$updates = @{
    1  = @("728÷5=", "867÷8=", "243÷4=", "331÷2=", "731÷8=")
    5  = @("737÷6=", "186÷6=", "866÷4=", "648÷6=", "542÷7=")
    9  = @("381÷9=", "551÷9=", "641÷7=", "939÷8=", "319÷9=")
    13 = @("927÷9=", "834÷8=", "749÷2=", "570÷9=", "845÷7=")
    17 = @("153÷7=", "651÷8=", "592÷2=", "124÷7=", "709÷8=")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    for ($col = 1; $col -le $values.Count; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$col - 1]
    }
}
